# Auto-generated script to update the cryptos worksheet with refreshed
# price / volume figures, matching the GitHub Actions data-refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.099.92'
$ws.Range('E2').Value = '  -3.46%  '
$ws.Range('D3').Value = '3.519.31'
$ws.Range('E3').Value = '  -4.59%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '581.26'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.29%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '175.38'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.35%  '
$ws.Range('E7').Value = '  +0.37%  '
$ws.Range('D8').Value = '3.512.47'
$ws.Range('E8').Value = '  -4.60%  '
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('E10').Value = '  -5.44%  '
$ws.Range('E11').Value = '  +7.55%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.602'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.77%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '47.41'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -4.87%  '
$ws.Range('E14').Value = '  -2.88%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '673.85'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -1.43%  '
$ws.Range('D16').Value = '4.089.04'
$ws.Range('E16').Value = '  -4.57%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '8.82'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.92%  '
$ws.Range('D18').Value = '3.523.74'
$ws.Range('E18').Value = '  -4.57%  '
$ws.Range('D19').Value = '69.117.92'
$ws.Range('E19').Value = '  -3.59%  '
$ws.Range('E20').Value = '  -1.45%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '17.62'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -2.56%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '11.28'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -3.23%  '
$ws.Range('E23').Value = '  -3.21%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '16.31'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -8.40%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '98.31'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -5.41%  '
$ws.Range('E26').Value = '  -4.13%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '5.85'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.19%  '
$ws.Range('E28').Value = '  -5.74%  '
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '9.51'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -6.55%  '
$ws.Range('E31').Value = '  -6.70%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '8.77'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -4.88%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.22'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -7.22%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '7.42'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +1.04%  '
$ws.Range('E35').Value = '  -5.16%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '578.97'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.91%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.63'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -13.80%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '10.96'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -3.04%  '
$ws.Range('E39').Value = '  -3.33%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '57.37'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -3.42%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.999'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('E42').Value = '  -3.39%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0442'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -4.06%  '
$ws.Range('E44').Value = '  -6.01%  '
$ws.Range('D45').Value = '3.438.83'
$ws.Range('E45').Value = '  -8.79%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '33.61'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -4.94%  '
$ws.Range('D47').Value = '0.0₃0710'
$ws.Range('E47').Value = '  -8.35%  '
$ws.Range('E48').Value = '  +1.16%  '
$ws.Range('E49').Value = '  -6.40%  '
$ws.Range('E50').Value = '  -0.23%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '131.55'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -2.17%  '
